$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.621.93"
$ws.Range("D2").Style = $s
$ws.Range("E2").Value = "  +0.44%  "

$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.928.16"
$ws.Range("D3").Style = $s
$ws.Range("E3").Value = "  +0.86%  "

$s = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("D4").Style = $s
$ws.Range("E4").Value = "  +0.43%  "

$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.51"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  +0.41%  "

$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4824"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  +0.24%  "

$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4065"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = "  +0.12%  "

$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08245"
$ws.Range("D9").Style = $s

$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.013"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  +0.29%  "

$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.94"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  +2.13%  "

$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.930.60"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  +0.84%  "

$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.102"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  +2.01%  "

$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.272"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = "  +2.20%  "

$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.92"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = "  +2.02%  "

$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06883"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  +1.66%  "

$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.012"
$ws.Range("D17").Style = $s
$ws.Range("E17").Value = "  +0.41%  "

$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = "  +0.38%  "

$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.66"
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("E20").Value = "  +0.46%  "

$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.604.39"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  +0.31%  "

$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.686"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  +1.29%  "

$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.00"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  +2.30%  "

$ws.Range("E24").Value = "  +0.42%  "

$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.167.64"
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = "  +0.98%  "

$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.27"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  +0.66%  "

$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.459"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  +1.45%  "

$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.07"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  +0.15%  "

$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.097"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  -0.14%  "

$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.89"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  +0.88%  "

$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.015"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  -1.12%  "

$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09645"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  +1.29%  "

$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.637"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  +2.53%  "

$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.579"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  +0.34%  "

$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.381"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = "  -0.32%  "

$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06406"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  +5.07%  "

$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02301"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  +1.58%  "

$ws.Range("E38").Value = "  +1.16%  "

$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5962"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  +0.53%  "

$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.78"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  +0.94%  "

$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.894"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  -0.66%  "

$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1854"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  +0.16%  "

$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.449"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  +1.51%  "

$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.286"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  +0.04%  "

$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.39"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  -0.93%  "

$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07551"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -1.16%  "

$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5574"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  +0.06%  "

$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.959"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  +1.22%  "

$s = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.55"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = "  +2.79%  "

$s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.442"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = "  +3.47%  "

$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.27"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -0.10%  "
